# Fill in a data row for the ProxyServer config sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string table insertion order: 127.0.0.1 (F2)
# must land first, then ProxyServer_1 (A2, reused by C2), then 000105001 (B2).
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "ProxyServer_1"
$ws.Range("B2").Value = "000105001"
$ws.Range("C2").Value = "ProxyServer_1"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 5001

# Move the active selection to G4 (no multi-cell selection)
$ws.Range("G4").Select()
